$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Proiect (column D) values for rows 9, 12, 18
$ws.Range("D9").Value = 2
$ws.Range("D12").Value = 2
$ws.Range("D18").Value = 2.25

# Recalculate the workbook so the Total (column F) formulas refresh
$excel.CalculateFullRebuild()

# Update the active cell / selection to D13
$ws.Range("D13").Select()
